$d = $word.ActiveDocument

# 1. Drop the "[Python, Javascript]" qualifier from the subtitle under the name.
$d.Content.Find.Execute(
    "Software & Web Developper [Python, Javascript]",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Software & Web Developper", 2
) | Out-Null

# 2. Add a new "Wordpress, SEO" bold line right after "Editor of choice: VIM"
#    in the Skills table cell.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Editor of choice:* VIM*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphAfter()
    $newIndex = $targetIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = "Wordpress, SEO"
    $newPara.Range.Font.Bold = $true
}
